# Update "想去人数" (column F) counts on the "展览" sheet and the
# corresponding rows on the "全部类型" sheet, per the commit's refreshed
# scrape data (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId=1): column F value updates by row.
$exhibitionUpdates = @{
    2  = 3140
    3  = 532
    4  = 1098
    5  = 87
    9  = 1131
    10 = 15786
    11 = 245
    12 = 185
    13 = 1028
    14 = 6196
    16 = 109
    28 = 31
    29 = 5004
    31 = 11089
    34 = 125
    35 = 175
    36 = 3808
    37 = 266
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" (sheetId=4): same updated counts, rows shifted because
# this sheet aggregates rows from the other category sheets.
$allTypesUpdates = @{
    3  = 3140
    4  = 532
    5  = 1098
    6  = 87
    10 = 1131
    11 = 15786
    12 = 245
    13 = 185
    14 = 1028
    15 = 6196
    17 = 109
    29 = 31
    30 = 5004
    33 = 11089
    36 = 125
    37 = 175
    38 = 3808
    39 = 266
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
